# Update cryptocurrency Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.410.23'
$ws.Range('E2').Value = '  -1.88%  '
# Row 3
$ws.Range('D3').Value = '3.011.62'
$ws.Range('E3').Value = '  -1.86%  '
# Row 4
$ws.Range('E4').Value = '  +0.00%  '
# Row 5
$ws.Range('D5').Value = "'584.47"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.51%  '
# Row 6
$ws.Range('D6').Value = "'147.88"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.05%  '
# Row 7
$ws.Range('E7').Value = '  +0.00%  '
# Row 8
$ws.Range('D8').Value = "'0.520"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.36%  '
# Row 9
$ws.Range('D9').Value = '3.003.47'
$ws.Range('E9').Value = '  -2.08%  '
# Row 10
$ws.Range('E10').Value = '  -4.13%  '
# Row 11
$ws.Range('D11').Value = "'5.70"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.19%  '
# Row 12
$ws.Range('D12').Value = "'0.442"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.23%  '
# Row 13
$ws.Range('D13').Value = "'0.0000230"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.80%  '
# Row 14
$ws.Range('D14').Value = "'34.76"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.29%  '
# Row 15
$ws.Range('E15').Value = '  +2.30%  '
# Row 16
$ws.Range('D16').Value = '3.499.20'
$ws.Range('E16').Value = '  -1.90%  '
# Row 17
$ws.Range('D17').Value = "'7.05"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.72%  '
# Row 18
$ws.Range('D18').Value = '62.299.66'
$ws.Range('E18').Value = '  -1.87%  '
# Row 19
$ws.Range('D19').Value = '3.005.49'
$ws.Range('E19').Value = '  -1.95%  '
# Row 20
$ws.Range('D20').Value = "'460.77"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.60%  '
# Row 21
$ws.Range('D21').Value = "'13.92"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.29%  '
# Row 22
$ws.Range('D22').Value = "'0.686"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.06%  '
# Row 23
$ws.Range('D23').Value = "'7.32"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.78%  '
# Row 24
$ws.Range('D24').Value = "'2.29"
$ws.Range('D24').Style = 'Normal'
# Row 25
$ws.Range('D25').Value = "'79.74"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.59%  '
# Row 26
$ws.Range('D26').Value = "'12.33"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.19%  '
# Row 27
$ws.Range('E27').Value = '  +0.11%  '
# Row 28
$ws.Range('D28').Value = "'9.96"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.22%  '
# Row 29
$ws.Range('E29').Value = '  +0.17%  '
# Row 30
$ws.Range('D30').Value = "'2.62"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.17%  '
# Row 31
$ws.Range('D31').Value = "'7.13"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.94%  '
# Row 32
$ws.Range('D32').Value = "'2.10"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.53%  '
# Row 33
$ws.Range('D33').Value = "'27.09"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.43%  '
# Row 34
$ws.Range('E34').Value = '  -4.09%  '
# Row 35
$ws.Range('E35').Value = '  -1.35%  '
# Row 36
$ws.Range('D36').Value = '0.0₃0788'
$ws.Range('E36').Value = '  -3.85%  '
# Row 37
$ws.Range('D37').Value = "'5.76"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.36%  '
# Row 38
$ws.Range('D38').Value = "'2.12"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.09%  '
# Row 39
$ws.Range('D39').Value = "'50.55"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.22%  '
# Row 40
$ws.Range('D40').Value = "'9.04"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.12%  '
# Row 41
$ws.Range('D41').Value = "'2.90"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -10.50%  '
# Row 42
$ws.Range('D42').Value = "'418.28"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.13%  '
# Row 43
$ws.Range('E43').Value = '  +1.21%  '
# Row 44
$ws.Range('D44').Value = "'0.275"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.50%  '
# Row 45
$ws.Range('D45').Value = '2.769.16'
$ws.Range('E45').Value = '  -0.87%  '
# Row 46
$ws.Range('D46').Value = "'0.0351"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.49%  '
# Row 47
$ws.Range('D47').Value = "'38.00"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.94%  '
# Row 48
$ws.Range('D48').Value = "'128.83"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.02%  '
# Row 49
$ws.Range('E49').Value = '  -0.02%  '
# Row 50
$ws.Range('E50').Value = '  -0.47%  '
# Row 51
$ws.Range('D51').Value = "'23.82"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.90%  '
